# Apply the change described by the diff:
#  - Productdata!G2: 40 -> 70
#  - ForecastedAverageDemand!B9:B11: 0 -> 100 (each)
#  - ForcastedStandardDeviation!B9:B11: 0 -> 10.23775 / 11.713975 / 13.0425775

$wb = $excel.ActiveWorkbook

$wsProductData = $wb.Worksheets.Item("Productdata")
$wsProductData.Range("G2").Value = 70

# The workbook contains a column (H) of cells that are "blank" but typed
# as strings (no shared-string value). A straight load/save round-trip of
# such cells through this runtime can turn them into the string at shared
# string index 0 ("Name"); explicitly re-blank them here so the untouched
# parts of the sheet stay identical to the source.
for ($r = 2; $r -le 11; $r++) {
    $wsProductData.Cells.Item($r, 8).Value = ""
}

$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
